$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 300 (new records for rows 301-302),
# pushing the former rows 301-349 down to 303-351.
$ws.Range("A301:A302").EntireRow.Insert()

# New row 301: Ají, Americana (o), Primera
$ws.Range("A301").Value = 2
$ws.Range("B301").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C301").Value = "Coquimbo"
$ws.Range("D301").Value = 44889
$ws.Range("E301").Value = 4
$ws.Range("F301").Value = 100112021
$ws.Range("G301").Value = "Ají"
$ws.Range("H301").Value = "Americana (o)"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 400
$ws.Range("K301").Value = 16000
$ws.Range("L301").Value = 18000
$ws.Range("M301").Value = 17000
$ws.Range("N301").Value = "`$/caja 25 kilos"
$ws.Range("O301").Value = "Provincia de Limarí"
$ws.Range("P301").Value = 680
$ws.Range("Q301").Value = 25
$ws.Range("R301").Value = "Hortaliza"

# New row 302: Ají, Inferno, Primera
$ws.Range("A302").Value = 2
$ws.Range("B302").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C302").Value = "Coquimbo"
$ws.Range("D302").Value = 44889
$ws.Range("E302").Value = 4
$ws.Range("F302").Value = 100112021
$ws.Range("G302").Value = "Ají"
$ws.Range("H302").Value = "Inferno"
$ws.Range("I302").Value = "Primera"
$ws.Range("J302").Value = 400
$ws.Range("K302").Value = 23000
$ws.Range("L302").Value = 25000
$ws.Range("M302").Value = 24000
$ws.Range("N302").Value = "`$/caja 25 kilos"
$ws.Range("O302").Value = "Provincia de Limarí"
$ws.Range("P302").Value = 960
$ws.Range("Q302").Value = 25
$ws.Range("R302").Value = "Hortaliza"
